$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the data rows (2-6) with placeholder field-name tokens,
# simulating a splice-in of database column names in place of literal
# row values. Header row (row 1) is left untouched.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = "sku"
    $ws.Cells.Item($r, 2).Value = "name"
    $ws.Cells.Item($r, 3).Value = "quantity"
    $ws.Cells.Item($r, 4).Value = "cost_per"
    $ws.Cells.Item($r, 5).Value = "total_cost"
}
